# Case and Fatality Demographics Data Updated
# Updates the weekly line-list numbers across all six sheets (three "Cases by …"
# sheets and three "Fatalities by …" sheets) and refreshes the workbook's
# selection / active-sheet UI state to match the authored edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Cases by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value  = 284
$ws.Range("B3").Value  = 1399
$ws.Range("B4").Value  = 3873
$ws.Range("B5").Value  = 15871
$ws.Range("B6").Value  = 17404
$ws.Range("B7").Value  = 15268
$ws.Range("B8").Value  = 12887
$ws.Range("B9").Value  = 4666
$ws.Range("B10").Value = 3155
$ws.Range("B11").Value = 1915
$ws.Range("B12").Value = 1267
$ws.Range("B13").Value = 1959
$ws.Activate()
$ws.Range("E10").Select()

# ---------------------------------------------------------------------------
# Cases by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 27316
$ws.Range("B3").Value = 51746

# ---------------------------------------------------------------------------
# Cases by RaceEthnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 979
$ws.Range("B3").Value = 13121
$ws.Range("B4").Value = 28632
$ws.Range("B5").Value = 581
$ws.Range("B6").Value = 27929
$ws.Range("B7").Value = 8719
$ws.Activate()
$ws.Range("D12").Select()

# ---------------------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B5").Value  = 269
$ws.Range("B6").Value  = 890
$ws.Range("B7").Value  = 2586
$ws.Range("B8").Value  = 5850
$ws.Range("B9").Value  = 4830
$ws.Range("B10").Value = 6215
$ws.Range("B11").Value = 6846
$ws.Range("B12").Value = 6739
$ws.Range("B13").Value = 16884
$ws.Activate()
$ws.Range("E7").Select()

# ---------------------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 21459
$ws.Range("B3").Value = 29703
$ws.Activate()
$ws.Range("D15").Select()

# ---------------------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1090
$ws.Range("B3").Value = 5223
$ws.Range("B4").Value = 23753
$ws.Range("B6").Value = 20794
$ws.Range("D13").Select()

# ---------------------------------------------------------------------------
# Final active sheet: "Fatalities by Gender" (4th zero-based tab) is the tab
# left active/selected in the saved workbook.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Fatalities by Gender").Activate()
